$d = $word.ActiveDocument

# 1) Merge the split runs around ", TractNHOPI, " into a single run
#    (removing the spell-check proofErr wrapping around "TractNHOPI").
$d.Content.Find.Execute(", TractNHOPI, ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ", TractNHOPI, ", 2) | Out-Null

# 2) Append two new empty paragraphs at the very end of the document body
#    (after the final "When analyzing food deserts..." paragraph, before
#    the sectPr). Using InsertXML with a minimal single-paragraph package
#    keeps the new paragraphs free of any placeholder run, matching how
#    Word represents a truly empty paragraph.
$emptyParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p></w:body></w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

for ($i = 0; $i -lt 2; $i++) {
    $endPos = $d.Content.End
    $endRange = $d.Range($endPos, $endPos)
    $endRange.InsertXML($emptyParaXml) | Out-Null
}

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
